# Commit: "updated data and pom file"
#
# The "Sheet1" worksheet holds two generated-test-user blocks:
#   rows 2..21  -> CONCATENATE($G,$I$2)       numeric suffix in I2  (14 -> 15)
#   rows 23..42 -> CONCATENATE($G,$I$23)      numeric suffix in I23 (19 -> 20)
# Columns A/B/C are formulas, so bumping the two suffix cells ripples
# through the whole block automatically.
#
# The "login" and "order" sheets hold a *static* (copy/pasted, not
# formula-linked) snapshot of that same generated data - login!G:I for
# rows 2-21 (the "14" block) and order!R:T for rows 2-21 (the "19"
# block). Those need the same suffix bump applied by hand since they are
# literal values, not formulas.

$wb = $excel.ActiveWorkbook

# ---- Sheet1: bump the two suffix cells that drive the CONCATENATE() formulas
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Range("I2").Value = 15
$sheet1.Range("I23").Value = 20

# ---- login sheet: static copies of the "14" block (rows 2-21), now "15"
$names14 = @(
    "EthanBaker", "DelanieCarman", "BretAgnew", "EdgardoTaylor", "TyrekReis",
    "LeannaChow", "TuckerCarlson", "AnnmarieConnor", "MoniqueWitte", "MikelWhitlock",
    "VincentAmaya", "KeiraQuiroz", "EllisCreech", "DionteCreel", "NicholeFoust",
    "ManuelConnell", "LourdesElam", "LincolnFrederick", "AlisaCash", "LucilleGriffiths"
)

$login = $wb.Worksheets.Item("login")
for ($i = 0; $i -lt $names14.Count; $i++) {
    $row = 2 + $i
    $user = $names14[$i] + "15"
    $email = $user + "@gmail.com"
    $login.Cells.Item($row, 7).Value = $user    # G
    $login.Cells.Item($row, 8).Value = $user    # H
    $login.Cells.Item($row, 9).Value = $email   # I
}

# ---- order sheet: static copies of the "19" block (rows 2-21), now "20"
$names19 = @(
    "DonnellJernigan", "MalikOtoole", "AlanCaudill", "AdanApplegate", "AiyanaWhitworth",
    "MercedezBrien", "DuaneHager", "LorenBell", "GeraldHiller", "DeionBranch",
    "DakotaHalstead", "ElliottFurman", "MiltonCamp", "DawnChester", "ZacheryPetrie",
    "EstebanAngel", "JimmyBlankenship", "AllysaGrice", "AugustineYoo", "BrandiSouthard"
)

$order = $wb.Worksheets.Item("order")
for ($i = 0; $i -lt $names19.Count; $i++) {
    $row = 2 + $i
    $user = $names19[$i] + "20"
    $email = $user + "@gmail.com"
    $order.Cells.Item($row, 18).Value = $user   # R
    $order.Cells.Item($row, 19).Value = $user   # S
    $order.Cells.Item($row, 20).Value = $email  # T
}

$excel.CalculateFullRebuild()
